# Replace English paragraph text with the Russian translations described
# in the diff.
#
# We look up each target paragraph by its current (English) text and then
# assign the new (Russian) text straight to that paragraph's Range. This
# keeps the edit anchored to the correct paragraph (unlike repeatedly
# re-deriving a fresh $d.Content range, whose position after Find.Execute
# is not reliably retained across statements) and preserves the existing
# xml:space="preserve" attribute on the w:t element, which a plain
# Find.Execute(...,replace,2) call strips whenever the replacement text has
# no leading/trailing whitespace.
$d = $word.ActiveDocument

function Replace-ParaText($find, $replace) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs($i)
        # Range.Text on a whole paragraph includes the trailing paragraph
        # mark (chr 13); strip it before comparing against $find.
        $ptext = $para.Range.Text
        if ($ptext.TrimEnd([char]13) -eq $find) {
            $para.Range.Text = $replace
            return
        }
    }
    throw "Paragraph text not found: $find"
}

Replace-ParaText "To confirm your registration, we would require you and one guest of your choice to provide us with:" "Чтобы подтвердить вашу регистрацию, нам необходимо, чтобы вы и ваш гость предоставили нам:"

Replace-ParaText "A scanned copy of your international passports" "Отсканированная копия ваших загранпаспортов"

Replace-ParaText "Covid-19 vaccination certificates" "Сертификаты о прививках Covid-19"

Replace-ParaText "Your country manager will be in touch to confirm your booking or request any other relevant details. " "Региональный менеджер свяжется с вами, чтобы подтвердить бронирование или запросить другие необходимые детали. "

Replace-ParaText "Our event package offers you and your guest: " "Наш пакет для мероприятия предлагает вам и вашему гостю: "

Replace-ParaText "Flight tickets " "Авиабилеты "

Replace-ParaText "Travel insurance " "Страхование путешествий "

Replace-ParaText "Airport – Hotel – Airport transfer " "Трансфер аэропорт – отель – аэропорт "

Replace-ParaText "One hotel room for you and your guest / Two hotel rooms for you and your guest" "Один гостиничный номер для вас и гостя / Два гостиничных номера для вас и гостя"

Replace-ParaText "Meals (Breakfast, lunch, and dinner)" "Питание (завтрак, обед и ужин)"

Replace-ParaText "We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. " "Перед датой вашего отъезда мы вышлем вам письмо-подтверждение с программой мероприятия и информацией о перелетах, транспорте и проживании. "

Replace-ParaText "We look forward to seeing you soon." "Мы будем рады встретиться с вами."
